$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header text in H1: "Tahun Pembelian" -> "Tanggal Pembelian (DD-MM-YYYY)"
$ws.Range("H1").Value = "Tanggal Pembelian (DD-MM-YYYY)"

# Widen column H to fit the new, longer header text (closest value this
# engine's character-width quantization can represent to the target 31.7109375)
$ws.Columns.Item(8).ColumnWidth = 30.8333333333333

# Scroll the view so column D is the left-most visible column, and select H5
$excel.ActiveWindow.ScrollColumn = 4
$ws.Range("H5").Select()
